$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $s = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $s
}

$ws.Range("D2").Value = '41.825.64'
$ws.Range("E2").Value = '  +4.23%  '

$ws.Range("D3").Value = '2.266.63'
$ws.Range("E3").Value = '  +2.01%  '

$ws.Range("E4").Value = '  +0.05%  '

Set-TextValue "D5" '305.11'
$ws.Range("E5").Value = '  +3.76%  '

Set-TextValue "D6" '92.01'
$ws.Range("E6").Value = '  +4.70%  '

$ws.Range("E7").Value = '  +3.75%  '

$ws.Range("E8").Value = '  +0.00%  '

Set-TextValue "D9" '0.483'
$ws.Range("E9").Value = '  +2.82%  '

Set-TextValue "D10" '32.65'
$ws.Range("E10").Value = '  +6.62%  '

Set-TextValue "D11" '53.88'
$ws.Range("E11").Value = '  +5.91%  '

$ws.Range("E12").Value = '  +2.06%  '

Set-TextValue "D14" '6.61'
$ws.Range("E14").Value = '  +3.23%  '

$ws.Range("D15").Value = '2.617.84'
$ws.Range("E15").Value = '  +1.98%  '

Set-TextValue "D16" '14.23'
$ws.Range("E16").Value = '  +2.96%  '

$ws.Range("D17").Value = '2.271.01'
$ws.Range("E17").Value = '  +2.05%  '

Set-TextValue "D18" '0.764'
$ws.Range("E18").Value = '  +3.87%  '

$ws.Range("D19").Value = '41.751.68'
$ws.Range("E19").Value = '  +4.22%  '

Set-TextValue "D20" '12.22'
$ws.Range("E20").Value = '  +8.53%  '

$ws.Range("D21").Value = '0.0₃0907'
$ws.Range("E21").Value = '  +1.99%  '

Set-TextValue "D22" '5.92'
$ws.Range("E22").Value = '  +2.40%  '

Set-TextValue "D23" '66.88'
$ws.Range("E23").Value = '  +1.87%  '

Set-TextValue "D24" '241.30'
$ws.Range("E24").Value = '  +2.39%  '

$ws.Range("E25").Value = '  +4.99%  '

$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("E27").Value = '  +5.35%  '

Set-TextValue "D28" '24.16'
$ws.Range("E28").Value = '  +4.24%  '

$ws.Range("E29").Value = '  +11.48%  '

Set-TextValue "D30" '9.62'
$ws.Range("E30").Value = '  +2.99%  '

Set-TextValue "D31" '34.04'
$ws.Range("E31").Value = '  +6.97%  '

Set-TextValue "D32" '157.68'
$ws.Range("E32").Value = '  -1.20%  '

$ws.Range("E33").Value = '  +0.03%  '

Set-TextValue "D34" '5.16'
$ws.Range("E34").Value = '  +3.97%  '

Set-TextValue "D35" '0.0747'
$ws.Range("E35").Value = '  +4.56%  '

Set-TextValue "D36" '3.02'
$ws.Range("E36").Value = '  -0.68%  '

Set-TextValue "D37" '17.25'
$ws.Range("E37").Value = '  +10.44%  '

$ws.Range("E38").Value = '  +1.47%  '

Set-TextValue "D39" '0.116'
$ws.Range("E39").Value = '  +2.75%  '

$ws.Range("E40").Value = '  +3.90%  '

$ws.Range("E41").Value = '  +2.33%  '

$ws.Range("E42").Value = '  +4.43%  '

$ws.Range("D43").Value = '2.063.39'
$ws.Range("E43").Value = '  -0.56%  '

Set-TextValue "D44" '19.40'
$ws.Range("E44").Value = '  -0.90%  '

$ws.Range("E45").Value = '  +3.06%  '

Set-TextValue "D46" '10.28'
$ws.Range("E46").Value = '  +2.83%  '

Set-TextValue "D47" '2.90'
$ws.Range("E47").Value = '  +4.68%  '

$ws.Range("E48").Value = '  +6.81%  '

$ws.Range("E49").Value = '  +4.52%  '

$ws.Range("E50").Value = '  +2.43%  '

Set-TextValue "D51" '72.50'
$ws.Range("E51").Value = '  +6.92%  '
